$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (non-numeric) cell updates: Coin name & Link swap for rows 41/42 ---
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

# --- Numeric-looking text cells (Price / Volume columns) ---
# Force Text format first so Excel keeps these as literal strings instead of
# auto-converting them to numbers/percentages (matches original inlineStr cells).
$numericCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "D12",
    "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17",
    "E17", "D18", "E18", "E20", "D21", "E21", "D22", "E22", "D23", "E23",
    "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D40",
    "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45",
    "E45", "E46", "D47", "E47", "D48", "E48", "E49", "E50"
)
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$values = @{
    "D2" = "256.59"
    "E2" = "0.58%"
    "D3" = "27.11"
    "E3" = "-1.41%"
    "D4" = "4.602"
    "E4" = "-11.29%"
    "D5" = "0.05910"
    "E5" = "0.93%"
    "D6" = "6.641"
    "E6" = "-1.30%"
    "D7" = "0.8694"
    "E7" = "0.06%"
    "D8" = "0.9438"
    "E8" = "-2.62%"
    "E9" = "0.03%"
    "D10" = "0.03750"
    "E10" = "9.14%"
    "D11" = "0.07086"
    "E11" = "-1.07%"
    "D12" = "0.03202"
    "E12" = "-0.32%"
    "D13" = "0.09261"
    "E13" = "0.39%"
    "D14" = "0.001546"
    "E14" = "-0.34%"
    "D15" = "0.0006046"
    "E15" = "-94.28%"
    "D16" = "0.006061"
    "E16" = "0.43%"
    "D17" = "3.512"
    "E17" = "0.41%"
    "D18" = "3.195"
    "E18" = "-1.17%"
    "E20" = "-1.48%"
    "D21" = "0.1284"
    "E21" = "0.28%"
    "D22" = "3.846"
    "E22" = "7.81%"
    "D23" = "0.04236"
    "E23" = "1.12%"
    "E24" = "-1.34%"
    "D25" = "0.001220"
    "E25" = "-0.12%"
    "D26" = "0.004290"
    "E26" = "-10.42%"
    "D27" = "0.0001200"
    "E27" = "0.01%"
    "D28" = "0.0001503"
    "E28" = "2.52%"
    "D40" = "0.03820"
    "E40" = "0.14%"
    "D41" = "0.006200"
    "E41" = "9.55%"
    "D42" = "0.1099"
    "E42" = "-0.16%"
    "D43" = "0.002277"
    "E43" = "-0.97%"
    "D44" = "0.01145"
    "E44" = "16.55%"
    "D45" = "0.00005506"
    "E45" = "2.28%"
    "E46" = "0.01%"
    "D47" = "0.08058"
    "E47" = "-19.41%"
    "D48" = "0.002281"
    "E48" = "7.17%"
    "E49" = "0.01%"
    "E50" = "0.01%"
}
foreach ($addr in $numericCells) {
    $ws.Range($addr).Value = $values[$addr]
}
